# Subindo esboço da aula 04
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly-tracked "check" dates for rows 10-14 (column C/D) ---
# Row 10: only the check date (D10) is new; C10 already has a value.
$ws.Range("C10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("D10").Value = 45381

# Rows 11-14: both start (C) and check (D) dates are new.
foreach ($r in 11..14) {
    $ws.Range("C6").Copy() | Out-Null
    $ws.Range("C$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("D$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("C$r").Value = 45381
    $ws.Range("D$r").Value = 45381
}

# Row 15: new row, only the date columns are populated.
$ws.Range("C6").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = 45384
$ws.Range("D15").Value = 45384

$excel.CutCopyMode = 0

# Row 20: new empty cell carrying a fresh underline style (no fill-in value).
$ws.Range("D20").Font.Underline = $true

# Selection follows the last-used cell, like Excel records on save.
$ws.Range("D20").Select() | Out-Null

# Print setup for the worksheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
